# Fruta / hortaliza, semanal
# Insert a new week's worth of data (2 rows) at the top of the date-sorted
# "Platano" price block (rows 626-627), pushing the existing data down by
# two rows. The sheet's used range grows from A1:T714 to A1:T716.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 626, shifting everything below (626-714) down
# to (628-716).
$ws.Rows("626:627").Insert()

# New row 626: "Pintón" quality record for the new week (Fecha 44776).
$ws.Range("A626").Value2 = 5
$ws.Range("B626").Value2 = "Macroferia Regional de Talca"
$ws.Range("C626").Value2 = "Maule"
$ws.Range("D626").Value2 = 44776
$ws.Range("E626").Value2 = 7
$ws.Range("F626").Value2 = "Fruta"
$ws.Range("G626").Value2 = 100108
$ws.Range("H626").Value2 = "Tropicales y subtropicales"
$ws.Range("I626").Value2 = 100108006
$ws.Range("J626").Value2 = "Plátano"
$ws.Range("K626").Value2 = "Sin especificar"
$ws.Range("L626").Value2 = "Pintón"
$ws.Range("M626").Value2 = 850
$ws.Range("N626").Value2 = 24000
$ws.Range("O626").Value2 = 24000
$ws.Range("P626").Value2 = 24000
$ws.Range("Q626").Value2 = "`$/caja 20 kilos"
$ws.Range("R626").Value2 = "Ecuador"
$ws.Range("S626").Value2 = 1200
$ws.Range("T626").Value2 = 20

# New row 627: "Primera Pintón" quality record for the same new week.
$ws.Range("A627").Value2 = 5
$ws.Range("B627").Value2 = "Macroferia Regional de Talca"
$ws.Range("C627").Value2 = "Maule"
$ws.Range("D627").Value2 = 44776
$ws.Range("E627").Value2 = 7
$ws.Range("F627").Value2 = "Fruta"
$ws.Range("G627").Value2 = 100108
$ws.Range("H627").Value2 = "Tropicales y subtropicales"
$ws.Range("I627").Value2 = 100108006
$ws.Range("J627").Value2 = "Plátano"
$ws.Range("K627").Value2 = "Sin especificar"
$ws.Range("L627").Value2 = "Primera Pintón"
$ws.Range("M627").Value2 = 600
$ws.Range("N627").Value2 = 25000
$ws.Range("O627").Value2 = 25000
$ws.Range("P627").Value2 = 25000
$ws.Range("Q627").Value2 = "`$/caja 20 kilos"
$ws.Range("R627").Value2 = "Ecuador"
$ws.Range("S627").Value2 = 1250
$ws.Range("T627").Value2 = 20
